$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (instance 1) - status remains OPTIMAL
$ws.Range("B2").Value = -274.9655845122987
$ws.Range("C2").Value = 0.09045029071781731
$ws.Range("D2").Value = 3418.929330501

# Row 3 (instance 2) - status changes to TIME_LIMIT
$ws.Range("B3").Value = -273.98162095885345
$ws.Range("C3").Value = 3.368418788147876
$ws.Range("D3").Value = 3633.83075055
$ws.Range("E3").Value = "TIME_LIMIT"

# Row 4 (instance 3) - status changes to TIME_LIMIT
$ws.Range("B4").Value = -274.0896045963642
$ws.Range("C4").Value = 3.2927143918754673
$ws.Range("D4").Value = 3671.370814071
$ws.Range("E4").Value = "TIME_LIMIT"

# Row 5 (instance 4) - status changes to TIME_LIMIT
$ws.Range("B5").Value = -276.8685496566217
$ws.Range("C5").Value = 4.176958512126726
$ws.Range("D5").Value = 3600.891001455
$ws.Range("E5").Value = "TIME_LIMIT"

# Row 6 (instance 5) - status changes to TIME_LIMIT
$ws.Range("B6").Value = -272.16539220117124
$ws.Range("C6").Value = 3.6729305785275232
$ws.Range("D6").Value = 3601.132246417
$ws.Range("E6").Value = "TIME_LIMIT"

# Row 7 (instance 6) - status remains OPTIMAL
$ws.Range("B7").Value = -268.97221193176233
$ws.Range("C7").Value = 0.014916014855304804
$ws.Range("D7").Value = 2156.109449285

# Row 8 (instance 7) - status changes to TIME_LIMIT
$ws.Range("B8").Value = -265.4281513734784
$ws.Range("C8").Value = 0.2929961414468023
$ws.Range("D8").Value = 3824.186122043
$ws.Range("E8").Value = "TIME_LIMIT"

# Row 9 (instance 8) - status changes to TIME_LIMIT
$ws.Range("B9").Value = -274.2017067884772
$ws.Range("C9").Value = 5.69475254826527
$ws.Range("D9").Value = 3638.063507639
$ws.Range("E9").Value = "TIME_LIMIT"

# Row 10 (instance 9) - status changes to TIME_LIMIT
$ws.Range("B10").Value = -271.53604103234676
$ws.Range("C10").Value = 0.22802400293768618
$ws.Range("D10").Value = 3816.246211745
$ws.Range("E10").Value = "TIME_LIMIT"

# Row 11 (instance 10) - status remains OPTIMAL
$ws.Range("B11").Value = -268.7867634966758
$ws.Range("C11").Value = 0.09885584246248551
$ws.Range("D11").Value = 2587.992979225
